$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados")
$ws.Activate()

# Populate columns D-I (valores de peças/serviços) for rows 138-185 that were
# added in this upload. Columns D-I already default to style index 11
# (numFmtId 2, "0.00") via the sheet's <cols> definition, so plain numeric
# assignment reproduces the same cell style as the diff.
$ws.Range("D138").Value = 1826.31
$ws.Range("E138").Value = 210.15
$ws.Range("D139").Value = 939.69
$ws.Range("E139").Value = 108.85
$ws.Range("D140").Value = 664.76
$ws.Range("E140").Value = 102
$ws.Range("G140").Value = 881.5
$ws.Range("D141").Value = 1380.29
$ws.Range("E141").Value = 69.8
$ws.Range("G141").Value = 284.08999999999997
$ws.Range("H141").Value = 73.8
$ws.Range("D142").Value = 820.98
$ws.Range("E142").Value = 114.95
$ws.Range("D143").Value = 273.42
$ws.Range("D144").Value = 275.48
$ws.Range("D145").Value = 227.15
$ws.Range("E145").Value = 154.30000000000001
$ws.Range("D146").Value = 152.6
$ws.Range("E146").Value = 154.30000000000001
$ws.Range("D147").Value = 359.43
$ws.Range("E147").Value = 183.1
$ws.Range("E148").Value = 37.79
$ws.Range("D149").Value = 786.01
$ws.Range("E149").Value = 177.13
$ws.Range("I149").Value = 228.58
$ws.Range("D150").Value = 1286.24
$ws.Range("E150").Value = 212.88
$ws.Range("I150").Value = 228.58
$ws.Range("E151").Value = 128.29
$ws.Range("G153").Value = 761.19
$ws.Range("I153").Value = 400.59
$ws.Range("E154").Value = 124.78
$ws.Range("D158").Value = 204.14
$ws.Range("D159").Value = 437.15
$ws.Range("E159").Value = 248.35
$ws.Range("G160").Value = 690.04
$ws.Range("H160").Value = 198.08
$ws.Range("E161").Value = 11.97
$ws.Range("E162").Value = 34.700000000000003
$ws.Range("D163").Value = 418.44
$ws.Range("E163").Value = 42.35
$ws.Range("G163").Value = 206.36
$ws.Range("D164").Value = 653.80999999999995
$ws.Range("E164").Value = 103.6
$ws.Range("G164").Value = 385.46
$ws.Range("D165").Value = 168.25
$ws.Range("D166").Value = 897.71
$ws.Range("E166").Value = 190.1
$ws.Range("G166").Value = 177.94
$ws.Range("D167").Value = 334.01
$ws.Range("D168").Value = 246.8
$ws.Range("G168").Value = 734.05
$ws.Range("D169").Value = 181.45
$ws.Range("E169").Value = 54.05
$ws.Range("D170").Value = 343.8
$ws.Range("E170").Value = 160.69999999999999
$ws.Range("D171").Value = 243.99
$ws.Range("E171").Value = 127
$ws.Range("E172").Value = 6.1
$ws.Range("D173").Value = 1241.47
$ws.Range("E173").Value = 240.1
$ws.Range("D174").Value = 732.56
$ws.Range("E174").Value = 125
$ws.Range("E175").Value = 81.849999999999994
$ws.Range("G176").Value = 143.80000000000001
$ws.Range("E178").Value = 113.65
$ws.Range("D181").Value = 1289.45
$ws.Range("D182").Value = 162.1
$ws.Range("D183").Value = 343.23
$ws.Range("E183").Value = 219.3
$ws.Range("F183").Value = 2.1
$ws.Range("E185").Value = 5.55

# Restore the scroll position / active-cell selection recorded for this sheet.
$ws.Range("H170").Select()
